# Remove the 'mrsow' and 'tsland' rows (and the two blank rows that followed
# them) from the "identified missing CMIP6 requested variables" sheet, so
# that these two variables become available again (issue #394).
#
# In the original workbook these entries occupy rows 70-71 (with rows 72-73
# being blank spacer rows before the next entry at row 74). Deleting the
# whole 70:73 block shifts every following row up by four, which also lets
# Excel drop the now-unused shared strings for "mrsow"/"tsland" and all of
# their associated metadata (long name, units, comment, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("70:73").Delete()
